$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: drop the placeholder empty cells (P3, R3, T3, U3) for this snapshot.
$ws.Range("P3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()

# New row 4: next classification snapshot for SH000016 (pre-market).
$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "SH000016"
$ws.Range("C4").Value = "上证50"
$ws.Range("D4").Value = "未开盘"
$ws.Range("E4").Value = "'2414.16"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "-0.75  -0.03%"
$ws.Range("G4").Value = 2420.98
$ws.Range("H4").Value = 2406.95
$ws.Range("I4").Value = 39706700
$ws.Range("J4").Value = 0.003
$ws.Range("K4").Value = 2403.27
$ws.Range("L4").Value = 2414.91
$ws.Range("M4").Value = 62284000000
$ws.Range("N4").Value = 1.26
$ws.Range("O4").Value = 2895.35
$ws.Range("P4").Value = "'"
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value = 0.007
$ws.Range("R4").Value = "'"
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").Value = 2288.01
$ws.Range("T4").Value = "'"
$ws.Range("T4").Style = "Normal"
$ws.Range("U4").Value = "'"
$ws.Range("U4").Style = "Normal"
